$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 172.5
$ws.Range("I39").Value = 172.5
$ws.Range("K39").Value = 517.5
$ws.Range("M39").Value = -221.5
$ws.Range("H42").Value = 1179.75
$ws.Range("I42").Value = 991.1429000000001
$ws.Range("J42").Value = 2500
$ws.Range("K42").Value = 2973.4287
$ws.Range("L42").Value = 7500
$ws.Range("M42").Value = -2743.4287
$ws.Range("N42").Value = -7960
$ws.Range("H70").Value = 2616.3333
$ws.Range("I70").Value = 2425
$ws.Range("J70").Value = 2999
$ws.Range("K70").Value = 7275
$ws.Range("L70").Value = 8997
$ws.Range("M70").Value = -7005
$ws.Range("N70").Value = -9537
$ws.Range("H73").Value = 2616.3333
$ws.Range("I73").Value = 2425
$ws.Range("J73").Value = 2999
$ws.Range("K73").Value = 7275
$ws.Range("L73").Value = 8997
$ws.Range("M73").Value = -6339
$ws.Range("N73").Value = -10869
$ws.Range("H76").Value = 66672332
$ws.Range("I76").Value = 83338340
$ws.Range("K76").Value = 83338340
$ws.Range("M76").Value = -83338025
$ws.Range("H79").Value = 66672332
$ws.Range("I79").Value = 83338340
$ws.Range("K79").Value = 83338340
$ws.Range("M79").Value = -83337248
$ws.Range("H82").Value = 8165
$ws.Range("I82").Value = 8165
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 24495
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -24089
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 8165
$ws.Range("I85").Value = 8165
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 24495
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -23091
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 2220.6191
$ws.Range("I86").Value = 2278
$ws.Range("J86").Value = 2105.8572
$ws.Range("K86").Value = 2278
$ws.Range("L86").Value = 2105.8572
$ws.Range("M86").Value = -1155
$ws.Range("N86").Value = -4351.8572
$ws.Range("H89").Value = 2220.6191
$ws.Range("I89").Value = 2278
$ws.Range("J89").Value = 2105.8572
$ws.Range("K89").Value = 11390
$ws.Range("L89").Value = 10529.286
$ws.Range("M89").Value = -5774
$ws.Range("N89").Value = -21761.286
$ws.Range("H101").Value = 1916.8
$ws.Range("J101").Value = 579
$ws.Range("L101").Value = 1737
$ws.Range("N101").Value = -4981
$ws.Range("H116").Value = 10104956
$ws.Range("I116").Value = 22225422
$ws.Range("J116").Value = 4568.3335
$ws.Range("K116").Value = 22225422
$ws.Range("L116").Value = 4568.3335
$ws.Range("M116").Value = -22221980
$ws.Range("N116").Value = -11452.3335
$ws.Range("H132").Value = 176924.83
$ws.Range("I132").Value = 223639.5
$ws.Range("K132").Value = 670918.5
$ws.Range("M132").Value = -668388.5
$ws.Range("H137").Value = 6283.3335
$ws.Range("I137").Value = 4167.147
$ws.Range("J137").Value = 11422.643
$ws.Range("K137").Value = 12501.441
$ws.Range("L137").Value = 34267.929
$ws.Range("M137").Value = -9951.440999999999
$ws.Range("N137").Value = -39367.929
$ws.Range("H138").Value = 4293.0444
$ws.Range("I138").Value = 494.3
$ws.Range("J138").Value = 5378.4
$ws.Range("K138").Value = 1482.9
$ws.Range("L138").Value = 16135.2
$ws.Range("M138").Value = 3657.1
$ws.Range("N138").Value = -26415.2
# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2800
$ws.Range("I45").Value = 2800
$ws.Range("K45").Value = 2800
$ws.Range("M45").Value = -2423
$ws.Range("H63").Value = 374.5
$ws.Range("I63").Value = 332.66666
$ws.Range("K63").Value = 332.66666
$ws.Range("M63").Value = 353.33334
$ws.Range("H66").Value = 374.5
$ws.Range("I66").Value = 332.66666
$ws.Range("K66").Value = 1663.3333
$ws.Range("M66").Value = 1768.6667
$ws.Range("H88").Value = 4456.0835
$ws.Range("J88").Value = 7356.9165
$ws.Range("L88").Value = 7356.9165
$ws.Range("N88").Value = -8168.9165
$ws.Range("H91").Value = 4456.0835
$ws.Range("J91").Value = 7356.9165
$ws.Range("L91").Value = 7356.9165
$ws.Range("N91").Value = -10164.9165
$ws.Range("H132").Value = 804917.2
$ws.Range("I132").Value = 948610.1
$ws.Range("J132").Value = 148035.14
$ws.Range("K132").Value = 2845830.3
$ws.Range("L132").Value = 444105.42
$ws.Range("M132").Value = -2843300.3
$ws.Range("N132").Value = -449165.42
# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 5520.5137
$ws.Range("I22").Value = 2704.7878
$ws.Range("J22").Value = 28750.25
$ws.Range("K22").Value = 2704.7878
$ws.Range("L22").Value = 28750.25
$ws.Range("M22").Value = -2531.7878
$ws.Range("N22").Value = -29096.25
$ws.Range("H134").Value = 4183016.2
$ws.Range("I134").Value = 7160313.5
$ws.Range("J134").Value = 14800
$ws.Range("K134").Value = 21480940.5
$ws.Range("L134").Value = 44400
$ws.Range("M134").Value = -21478405.5
$ws.Range("N134").Value = -49470
# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 41681.668
$ws.Range("J70").Value = 44018
$ws.Range("L70").Value = 44018
$ws.Range("N70").Value = -44648
$ws.Range("H73").Value = 41681.668
$ws.Range("J73").Value = 44018
$ws.Range("L73").Value = 44018
$ws.Range("N73").Value = -46202
$ws.Range("H132").Value = 4112.587
$ws.Range("I132").Value = 3503.611
$ws.Range("K132").Value = 10510.833
$ws.Range("M132").Value = -7980.832999999999
$ws.Range("H134").Value = 50007296
$ws.Range("I134").Value = 83340470
$ws.Range("K134").Value = 250021410
$ws.Range("M134").Value = -250018875
# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 96250.39999999999
$ws.Range("I9").Value = 459.5
$ws.Range("K9").Value = 1378.5
$ws.Range("M9").Value = -1154.5
$ws.Range("H113").Value = 4091980
$ws.Range("J113").Value = 501177.9
$ws.Range("L113").Value = 1503533.7
$ws.Range("N113").Value = -1507873.7
$ws.Range("H128").Value = 147665.83
$ws.Range("I128").Value = 147665.83
$ws.Range("K128").Value = 442997.49
$ws.Range("M128").Value = -438017.49
# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 632
$ws.Range("J107").Value = 664
$ws.Range("L107").Value = 664
$ws.Range("N107").Value = -4504
$ws.Range("H126").Value = 15629619
$ws.Range("J126").Value = 7558.273
$ws.Range("L126").Value = 22674.819
$ws.Range("N126").Value = -27614.819
$ws.Range("H132").Value = 20003518
$ws.Range("I132").Value = 32261324
$ws.Range("J132").Value = 3940.5264
$ws.Range("K132").Value = 96783972
$ws.Range("L132").Value = 11821.5792
$ws.Range("M132").Value = -96781442
$ws.Range("N132").Value = -16881.5792
# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4831.095
$ws.Range("J40").Value = 5808.625
$ws.Range("L40").Value = 5808.625
$ws.Range("N40").Value = -6080.625
$ws.Range("H68").Value = 1738.4286
$ws.Range("I68").Value = 1382.1765
$ws.Range("K68").Value = 1382.1765
$ws.Range("M68").Value = -633.1765
$ws.Range("H71").Value = 1738.4286
$ws.Range("I71").Value = 1382.1765
$ws.Range("K71").Value = 6910.8825
$ws.Range("M71").Value = -3166.8825
$ws.Range("H132").Value = 5020.587
$ws.Range("I132").Value = 4187.222
$ws.Range("K132").Value = 12561.666
$ws.Range("M132").Value = -10031.666
# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15794.889
$ws.Range("I62").Value = 11343.571
$ws.Range("J62").Value = 31374.5
$ws.Range("K62").Value = 11343.571
$ws.Range("L62").Value = 31374.5
$ws.Range("M62").Value = -10719.571
$ws.Range("N62").Value = -32622.5
$ws.Range("H65").Value = 15794.889
$ws.Range("I65").Value = 11343.571
$ws.Range("J65").Value = 31374.5
$ws.Range("K65").Value = 56717.855
$ws.Range("L65").Value = 156872.5
$ws.Range("M65").Value = -53597.855
$ws.Range("N65").Value = -163112.5
$ws.Range("H132").Value = 6086.8
$ws.Range("I132").Value = 5371.4287
$ws.Range("K132").Value = 16114.2861
$ws.Range("M132").Value = -13584.2861
